# Apply updated crypto price/volume data (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.870.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.112.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.52%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.67"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.77%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.105.43"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.24%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.08%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +10.44%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000244"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +4.85%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.124"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.626.81"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.51%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.110.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.39%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "62.837.58"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "452.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.735"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.54"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.65"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.83"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.75%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.48%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.28"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.56%  "

$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.86"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +11.73%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.112"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +11.49%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.14"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.04%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.05"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0804"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.44%  "

$ws.Range("B36").Value = "Stacks"
$ws.Range("C36").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.30"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.25%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.08"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.88"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.81%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.23%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "425.78"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.952.26"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0375"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.68%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +9.66%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +7.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.02"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.76%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.82"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.96%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.06%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.98%  "
